# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values for column G (rows 2-31), replacing the old Strike# values.
$kValues = @{
    2  = 2
    3  = 9
    4  = 1
    5  = 11
    6  = 6
    7  = 7
    8  = 10
    9  = 2
    10 = 5
    11 = 2
    12 = 6
    13 = 6
    14 = 5
    15 = 6
    16 = 8
    17 = 7
    18 = 8
    19 = 3
    20 = 4
    21 = 4
    22 = 5
    23 = 7
    24 = 6
    25 = 7
    26 = 6
    27 = 5
    28 = 4
    29 = 1
    30 = 4
    31 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
